# Auto-generated Excel COM-interop script to update the cryptos price table
# on Sheet1 to match the target diff (price/volume refresh + two row swaps).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '37.789.21'
$ws.Range('E2').Value = '  +0.79%  '

# Row 3
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.084.78'
$ws.Range('E3').Value = '  +0.20%  '

# Row 4
$ws.Range('E4').Value = '  +0.02%  '

# Row 5
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '234.90'
$ws.Range('E5').Value = '  +0.48%  '

# Row 6
$ws.Range('E6').Value = '  +0.27%  '

# Row 7
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '58.91'
$ws.Range('E7').Value = '  +2.00%  '

# Row 8
$ws.Range('E8').Value = '  +0.00%  '

# Row 9
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.392'
$ws.Range('E9').Value = '  +1.12%  '

# Row 10
$ws.Range('E10').Value = '  +1.57%  '

# Row 11
$ws.Range('E11').Value = '  +2.73%  '

# Row 12
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.392.58'
$ws.Range('E12').Value = '  +0.23%  '

# Row 13
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '14.80'
$ws.Range('E13').Value = '  +2.75%  '

# Row 14
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '21.17'
$ws.Range('E14').Value = '  +1.50%  '

# Row 15
$ws.Range('E15').Value = '  -0.83%  '

# Row 16
$ws.Range('E16').Value = '  +2.18%  '

# Row 17
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '2.040.95'
$ws.Range('E17').Value = '  -1.89%  '

# Row 18
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '37.699.54'
$ws.Range('E18').Value = '  +0.73%  '

# Row 19
$ws.Range('E19').Value = '  -0.60%  '

# Row 20
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '71.40'

# Row 21
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0835'
$ws.Range('E21').Value = '  +1.96%  '

# Row 22
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '229.01'
$ws.Range('E22').Value = '  +1.46%  '

# Row 23
$ws.Range('E23').Value = '  -0.04%  '

# Row 24
$ws.Range('E24').Value = '  -2.07%  '

# Row 25
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.43'
$ws.Range('E25').Value = '  +1.58%  '

# Row 26
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '169.27'
$ws.Range('E26').Value = '  +0.69%  '

# Row 27
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.139'
$ws.Range('E27').Value = '  +6.13%  '

# Row 28
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '9.02'
$ws.Range('E28').Value = '  +1.75%  '

# Row 29
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '19.58'
$ws.Range('E29').Value = '  +2.20%  '

# Row 30
$ws.Range('E30').Value = '  -1.86%  '

# Row 31
$ws.Range('E31').Value = '  +2.48%  '

# Row 32
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.70'
$ws.Range('E32').Value = '  +1.73%  '

# Row 33
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.0635'
$ws.Range('E33').Value = '  +2.81%  '

# Row 34
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '4.68'
$ws.Range('E34').Value = '  +2.27%  '

# Row 35
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '2.50'
$ws.Range('E35').Value = '  -1.49%  '

# Row 36
$ws.Range('E36').Value = '  +3.18%  '

# Row 37
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '3.39'
$ws.Range('E37').Value = '  -1.95%  '

# Row 38
$ws.Range('E38').Value = '  -0.19%  '

# Row 39
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.45'
$ws.Range('E39').Value = '  -1.90%  '

# Row 40
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.0981'
$ws.Range('E40').Value = '  +2.45%  '

# Row 41
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.83'
$ws.Range('E41').Value = '  +1.98%  '

# Row 42
$ws.Range('B42').Value = 'VeChain'
$ws.Range('C42').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0216'
$ws.Range('E42').Value = '  +1.50%  '

# Row 43
$ws.Range('B43').Value = 'HuobiToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '2.88'
$ws.Range('E43').Value = '  -2.08%  '

# Row 44
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.464.17'
$ws.Range('E44').Value = '  -1.32%  '

# Row 45
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '4.34'
$ws.Range('E45').Value = '  +4.55%  '

# Row 46
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '1.17'
$ws.Range('E46').Value = '  +1.04%  '

# Row 47
$ws.Range('B47').Value = 'InjectiveProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '16.28'
$ws.Range('E47').Value = '  +5.12%  '

# Row 48
$ws.Range('B48').Value = 'ARBITRUM'
$ws.Range('C48').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.07'
$ws.Range('E48').Value = '  +2.44%  '

# Row 49
$ws.Range('E49').Value = '  +3.16%  '

# Row 50
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '3.02'
$ws.Range('E50').Value = '  +0.84%  '

# Row 51
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '2.277.95'
$ws.Range('E51').Value = '  +0.22%  '
